$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Resposibility (column E) text values.
# Written in this order so the rebuilt shared-strings table lines up with
# the new text for each reused index (E7 first introduces the new string
# that ends up at index 19, then E2 -> 20, E3 -> 21, E4 -> 22).
$ws.Range("E7").Value = "Customer and Consumer Login and Registration"
$ws.Range("E2").Value = "Customer Level Operation"
$ws.Range("E3").Value = "Consumer Operations"
$ws.Range("E4").Value = "Bank Representative Level Account Operations"
$ws.Range("E5").Value = "Consumer Operations"
$ws.Range("E6").Value = "Customer Level Operation"

# Widen column E to fit the longer text (target stored width 47.44140625).
$ws.Columns.Item(5).ColumnWidth = 46.71

# Update the view: change the zoom level and move the active selection to
# E5 (this also clears the old scrolled-in topLeftCell="C1").
$excel.ActiveWindow.Zoom = 145
[void]$ws.Range("E5").Select()
